# PlayerPerformance_4851.xlsx edit
#  1. Remove the 5 empty INNING_NUMBER cells (B2,B5,B6,B7,B12) on "ODI Batting"
#  2. Add a new worksheet "ODI Batting Extra" (after "ODI Bowling") with
#     MATCH_CODE / BATTING_POSITION / NUM_4 / NUM_6 / PERCENT_RUNS_OF_TOTAL /
#     MAN_OF_MATCH columns for the 11 matches.

$wb = $excel.ActiveWorkbook

# --- 1. "ODI Batting": clear the empty INNING_NUMBER cells so they disappear ---
$odiBatting = $wb.Worksheets.Item("ODI Batting")
$odiBatting.Range("B2").ClearContents()
$odiBatting.Range("B5").ClearContents()
$odiBatting.Range("B6").ClearContents()
$odiBatting.Range("B7").ClearContents()
$odiBatting.Range("B12").ClearContents()

# --- 2. Add the new "ODI Batting Extra" sheet after "ODI Bowling" ---
$odiBowling = $wb.Worksheets.Item("ODI Bowling")
$newSheet = $wb.Worksheets.Add([Type]::Missing, $odiBowling, 1, [Type]::Missing)
$newSheet.Name = "ODI Batting Extra"

# Header row — reuse the same bold/border/center header style used by the
# other sheets in this workbook.
$odiBowling.Range("A1").Copy()
$newSheet.Range("A1:F1").PasteSpecial(-4122)  # xlPasteFormats

$newSheet.Range("A1").Value = "MATCH_CODE"
$newSheet.Range("B1").Value = "BATTING_POSITION"
$newSheet.Range("C1").Value = "NUM_4"
$newSheet.Range("D1").Value = "NUM_6"
$newSheet.Range("E1").Value = "PERCENT_RUNS_OF_TOTAL"
$newSheet.Range("F1").Value = "MAN_OF_MATCH"

# Helper: force a numeric-looking string to be stored as text (matches the
# source data, which keeps match codes / NUM_4 / NUM_6 / percentages as text).
function Set-TextValue($cell, [string]$val) {
    $cell.NumberFormat = "@"
    $cell.Value = $val
}

# Row 2 — match 4563
Set-TextValue $newSheet.Cells.Item(2,1) "4563"
$newSheet.Cells.Item(2,2).Value = 11
Set-TextValue $newSheet.Cells.Item(2,3) "0"
Set-TextValue $newSheet.Cells.Item(2,4) "0"
Set-TextValue $newSheet.Cells.Item(2,5) "0.38%"
$newSheet.Cells.Item(2,6).Value = "NO"

# Row 3 — match 4566
Set-TextValue $newSheet.Cells.Item(3,1) "4566"
$newSheet.Cells.Item(3,2).Value = 10
$newSheet.Cells.Item(3,6).Value = "NO"

# Row 4 — match 4605
Set-TextValue $newSheet.Cells.Item(4,1) "4605"
$newSheet.Cells.Item(4,6).Value = "NO"

# Row 5 — match 4608
Set-TextValue $newSheet.Cells.Item(5,1) "4608"
$newSheet.Cells.Item(5,2).Value = 11
Set-TextValue $newSheet.Cells.Item(5,3) "0"
Set-TextValue $newSheet.Cells.Item(5,4) "0"
Set-TextValue $newSheet.Cells.Item(5,5) "1.85%"
$newSheet.Cells.Item(5,6).Value = "NO"

# Row 6 — match 4614
Set-TextValue $newSheet.Cells.Item(6,1) "4614"
$newSheet.Cells.Item(6,2).Value = 11
Set-TextValue $newSheet.Cells.Item(6,3) "1"
Set-TextValue $newSheet.Cells.Item(6,4) "0"
Set-TextValue $newSheet.Cells.Item(6,5) "2.19%"
$newSheet.Cells.Item(6,6).Value = "NO"

# Row 7 — match 4625
Set-TextValue $newSheet.Cells.Item(7,1) "4625"
$newSheet.Cells.Item(7,6).Value = "NO"

# Row 8 — match 4692
Set-TextValue $newSheet.Cells.Item(8,1) "4692"

# Row 9 — match 4695
Set-TextValue $newSheet.Cells.Item(9,1) "4695"

# Row 10 — match 4697
Set-TextValue $newSheet.Cells.Item(10,1) "4697"

# Row 11 — match 4735
Set-TextValue $newSheet.Cells.Item(11,1) "4735"

# Row 12 — match 4745
Set-TextValue $newSheet.Cells.Item(12,1) "4745"

$newSheet.Range("A1").Select()
